$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.26762998866126964,
    -0.0099999998446342886,
    -0.0089999998430467798,
    0.061997998877032501,
    -0.0059999998472877181,
    -0.0059999998430093626,
    -0.019999999809321878,
    -0.01999999980782885,
    -0.0059999998395410259,
    -0.0059999998383233333,
    -0.0044999998418511211,
    -0.0059999998380924069,
    -0.0059999998373703178,
    0.020494951472663558,
    -0.0059999998366215834,
    -0.0059999998360180662,
    -0.0059999998352546768,
    -0.0089999998279477467,
    -0.0089999998468717202,
    -0.0089999998454235453,
    -0.0089999998451952834,
    -0.0089999998450052132,
    -0.048459575449909309,
    -0.041999999760057527,
    -0.041999999758678186,
    -0.0059999998427002765,
    -0.0059999998423077017,
    -0.0059999998407525013,
    0.03325763008875704,
    -0.019999999805444535,
    -0.014999999815573872,
    -0.020999999801046165,
    -0.0059999998364794749
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
